$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 17240.854
$ws.Range("I15").Value = 17240.854
$ws.Range("K15").Value = 51722.562
$ws.Range("M15").Value = -51553.562

$ws.Range("H28").Value = 522.4286
$ws.Range("I28").Value = 357.14285
$ws.Range("K28").Value = 357.14285
$ws.Range("M28").Value = 127.85715

$ws.Range("H132").Value = 1412.5555
$ws.Range("I132").Value = 1666.7273
$ws.Range("J132").Value = 294.2
$ws.Range("K132").Value = 5000.1819
$ws.Range("L132").Value = 882.5999999999999
$ws.Range("M132").Value = -2470.1819
$ws.Range("N132").Value = -5942.6

$ws.Range("H138").Value = 2977.2742
$ws.Range("I138").Value = 1207.0541
$ws.Range("J138").Value = 5597.2
$ws.Range("K138").Value = 3621.1623
$ws.Range("L138").Value = 16791.6
$ws.Range("M138").Value = 1518.8377
$ws.Range("N138").Value = -27071.6

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4624.8765
$ws.Range("I32").Value = 4017.7163
$ws.Range("J32").Value = 7530.5713
$ws.Range("K32").Value = 4017.7163
$ws.Range("L32").Value = 7530.5713
$ws.Range("M32").Value = -3730.7163
$ws.Range("N32").Value = -8104.5713

$ws.Range("H45").Value = 26335.334
$ws.Range("I45").Value = 38503
$ws.Range("J45").Value = 2000
$ws.Range("K45").Value = 38503
$ws.Range("L45").Value = 2000
$ws.Range("M45").Value = -38126
$ws.Range("N45").Value = -2754

$ws.Range("H74").Value = 3233.2693
$ws.Range("I74").Value = 3714.973
$ws.Range("J74").Value = 2045.0667
$ws.Range("K74").Value = 3714.973
$ws.Range("L74").Value = 2045.0667
$ws.Range("M74").Value = -2840.973
$ws.Range("N74").Value = -3793.0667

$ws.Range("H77").Value = 3233.2693
$ws.Range("I77").Value = 3714.973
$ws.Range("J77").Value = 2045.0667
$ws.Range("K77").Value = 18574.865
$ws.Range("L77").Value = 10225.3335
$ws.Range("M77").Value = -14206.865
$ws.Range("N77").Value = -18961.3335

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 981.2632
$ws.Range("I107").Value = 934.1
$ws.Range("K107").Value = 934.1
$ws.Range("M107").Value = 985.9

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 19231446
$ws.Range("I16").Value = 25641594
$ws.Range("J16").Value = 1000
$ws.Range("K16").Value = 25641594
$ws.Range("L16").Value = 1000
$ws.Range("M16").Value = -25641307
$ws.Range("N16").Value = -1574

$ws.Range("H105").Value = 37038584
$ws.Range("I105").Value = 55556456
$ws.Range("J105").Value = 2833.3333
$ws.Range("K105").Value = 55556456
$ws.Range("L105").Value = 2833.3333
$ws.Range("M105").Value = -55554709
$ws.Range("N105").Value = -6327.3333

$ws.Range("H107").Value = 559.1081
$ws.Range("I107").Value = 464.3158
$ws.Range("J107").Value = 659.1667
$ws.Range("K107").Value = 464.3158
$ws.Range("L107").Value = 659.1667
$ws.Range("M107").Value = 1455.6842
$ws.Range("N107").Value = -4499.1667

$ws.Range("H113").Value = 19231446
$ws.Range("I113").Value = 25641594
$ws.Range("J113").Value = 1000
$ws.Range("K113").Value = 25641594
$ws.Range("L113").Value = 1000
$ws.Range("M113").Value = -25639424
$ws.Range("N113").Value = -5340

$ws.Range("H132").Value = 2670.56
$ws.Range("I132").Value = 1566.125
$ws.Range("J132").Value = 4634
$ws.Range("K132").Value = 4698.375
$ws.Range("L132").Value = 13902
$ws.Range("M132").Value = -2168.375
$ws.Range("N132").Value = -18962

$ws.Range("H135").Value = 33791.668
$ws.Range("J135").Value = 33791.668
$ws.Range("L135").Value = 33791.668
$ws.Range("N135").Value = -43931.668

$ws.Range("H140").Value = 22137.057
$ws.Range("J140").Value = 22137.057
$ws.Range("L140").Value = 22137.057
$ws.Range("N140").Value = -32497.057

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 2895.4626
$ws.Range("I68").Value = 3068.721
$ws.Range("J68").Value = 2585.0417
$ws.Range("K68").Value = 9206.163
$ws.Range("L68").Value = 7755.125100000001
$ws.Range("M68").Value = -8395.163
$ws.Range("N68").Value = -9377.125100000001

$ws.Range("H71").Value = 2895.4626
$ws.Range("I71").Value = 3068.721
$ws.Range("J71").Value = 2585.0417
$ws.Range("K71").Value = 27618.489
$ws.Range("L71").Value = 23265.3753
$ws.Range("M71").Value = -23562.489
$ws.Range("N71").Value = -31377.3753

$ws.Range("H113").Value = 1304853.8
$ws.Range("I113").Value = 1852315.4
$ws.Range("J113").Value = 526882
$ws.Range("K113").Value = 5556946.199999999
$ws.Range("L113").Value = 1580646
$ws.Range("M113").Value = -5554776.199999999
$ws.Range("N113").Value = -1584986

$ws.Range("H120").Value = 8000
$ws.Range("J120").Value = 8000
$ws.Range("L120").Value = 24000
$ws.Range("N120").Value = -33676

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 10117.5
$ws.Range("I126").Value = 10117.5
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 30352.5
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -27882.5
$ws.Range("N126").ClearContents()

$ws.Range("H133").Value = 87090
$ws.Range("J133").Value = 87090
$ws.Range("L133").Value = 87090
$ws.Range("N133").Value = -97210

$ws.Range("H135").Value = 38222.855
$ws.Range("J135").Value = 38222.855
$ws.Range("L135").Value = 38222.855
$ws.Range("N135").Value = -48362.855

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("M4").ClearContents()

$ws.Range("H7").Value = 8839.077
$ws.Range("I7").Value = 2790.3
$ws.Range("K7").Value = 2790.3
$ws.Range("M7").Value = -2678.3

$ws.Range("H28").Value = 0
$ws.Range("I28").Value = 0
$ws.Range("K28").Value = 0
$ws.Range("M28").ClearContents()

$ws.Range("H37").Value = 0
$ws.Range("I37").Value = 0
$ws.Range("K37").Value = 0
$ws.Range("M37").ClearContents()

$ws.Range("H61").Value = 1689.8334
$ws.Range("I61").Value = 1689.8334
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 1689.8334
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -1487.8334
$ws.Range("N61").ClearContents()

$ws.Range("H113").Value = 1689.8334
$ws.Range("I113").Value = 1689.8334
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 1689.8334
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 480.1666
$ws.Range("N113").ClearContents()

$ws.Range("H126").Value = 8839.077
$ws.Range("I126").Value = 2790.3
$ws.Range("K126").Value = 8370.900000000001
$ws.Range("M126").Value = -5900.900000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2869.45
$ws.Range("I132").Value = 2942.3333
$ws.Range("J132").Value = 2760.125
$ws.Range("K132").Value = 8826.999899999999
$ws.Range("L132").Value = 8280.375
$ws.Range("M132").Value = -6296.999899999999
$ws.Range("N132").Value = -13340.375
